$d = $word.ActiveDocument

# --- Change 1 --------------------------------------------------------------
# "Viene cliccato il link "Hai dimenticato la password?", che porta alla
# pagina della password smarrita." used to be split across two runs with a
# "_GoBack" bookmark sandwiched between them. Re-merge it into a single run
# (removing that now-orphaned bookmark at the same time — it gets a new
# home further down the document as part of change 2).
$d.Content.Find.Execute(
    "Hai dimenticato la password", $true, $false, $false, $false, $false,
    $true, 1, $false, "Hai dimenticato la password", 2) | Out-Null

# --- Change 2 --------------------------------------------------------------
# Append the new "daily log" paragraphs after the last table, right before
# the final section break, describing today's work (and relocate the
# "_GoBack" bookmark there).

# Blank paragraph
$insertionPoint = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$p1 = $d.Paragraphs.Add($insertionPoint)

# "Oggi Thor inizialmente di " paragraph
$insertionPoint = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$p2 = $d.Paragraphs.Add($insertionPoint)
$p2.Range.Text = "Oggi Thor inizialmente di "

# Blank paragraph
$insertionPoint = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$p3 = $d.Paragraphs.Add($insertionPoint)

# Final, detailed paragraph
$insertionPoint = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$p4 = $d.Paragraphs.Add($insertionPoint)
$p4.Range.Text = "Oggi Thor inizialmente si è occupato di dare un’occhiata ai test di Selenium che riscontravano un errore la lezione precedente (riguardante la connessione e il display di xvfb), successivamente si è occupato dei Test Case, capitolo 4.1 della documentazione."

$d.Bookmarks.Add("_GoBack", $p4.Range) | Out-Null
